# Natmi following Dr Hou advice
# Updates LR-pair summary rows for ECs/FAPs/M2/sCs (Mif -> Cxcr4) and adds the sCs sending-cluster block (rows 14-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowsData = @(
    @{ A="ECs"; B="Mif"; C="Cxcr4"; D="ECs"; E=3; F=1; G=24.11008933333333; H=72.33026799999999; I=0.2144059681009565; J=0.2144059681009565; K=3; L=1; M=23.66111866666667; N=70.983356; O=0.08861204511346259; P=0.08861204511346259; Q=570.4716847799342; R=5134.245163019407; S=0.01899895131795758; T=0.01899895131795758 },
    @{ A="ECs"; B="Mif"; C="Cxcr4"; D="FAPs"; E=3; F=1; G=24.11008933333333; H=72.33026799999999; I=0.2144059681009565; J=0.2144059681009565; K=1; L=0.3333333333333333; M=0.1622346666666667; N=0.486704; O=0.0006075767508780888; P=0.0006075767508780888; Q=3.911492306296889; R=35.203430756672; S=0.0001302680814676503; T=0.0001302680814676503 },
    @{ A="ECs"; B="Mif"; C="Cxcr4"; D="M2"; E=3; F=1; G=24.11008933333333; H=72.33026799999999; I=0.2144059681009565; J=0.2144059681009565; K=3; L=1; M=231.2048796666667; N=693.614639; O=0.8658735673532548; P=0.8658735673532548; Q=5574.370303065916; R=50169.33272759325; S=0.1856484604614034; T=0.1856484604614033 },
    @{ A="ECs"; B="Mif"; C="Cxcr4"; D="sCs"; E=3; F=1; G=24.11008933333333; H=72.33026799999999; I=0.2144059681009565; J=0.2144059681009565; K=3; L=1; M=11.99098133333333; N=35.972944; O=0.04490681078240458; P=0.04490681078240458; Q=289.103631140999; R=2601.932680268992; S=0.009628288240127928; T=0.009628288240127925 },
    @{ A="FAPs"; B="Mif"; C="Cxcr4"; D="ECs"; E=3; F=1; G=27.64911833333333; H=82.947355; I=0.2458778107968398; J=0.2458778107968398; K=3; L=1; M=23.66111866666667; N=70.983356; O=0.08861204511346259; P=0.08861204511346259; Q=654.209069913709; R=5887.881629223381; S=0.02178773566272899; T=0.02178773566272899 },
    @{ A="FAPs"; B="Mif"; C="Cxcr4"; D="FAPs"; E=3; F=1; G=27.64911833333333; H=82.947355; I=0.2458778107968398; J=0.2458778107968398; K=1; L=0.3333333333333333; M=0.1622346666666667; N=0.486704; O=0.0006075767508780888; P=0.0006075767508780888; Q=4.485645496435556; R=40.37080946792; S=0.0001493896413969614; T=0.0001493896413969614 },
    @{ A="FAPs"; B="Mif"; C="Cxcr4"; D="M2"; E=3; F=1; G=27.64911833333333; H=82.947355; I=0.2458778107968398; J=0.2458778107968398; K=3; L=1; M=231.2048796666667; N=693.614639; O=0.8658735673532548; P=0.8658735673532548; Q=6392.611077147761; R=57533.49969432985; S=0.2128990971676683; T=0.2128990971676683 },
    @{ A="FAPs"; B="Mif"; C="Cxcr4"; D="sCs"; E=3; F=1; G=27.64911833333333; H=82.947355; I=0.2458778107968398; J=0.2458778107968398; K=3; L=1; M=11.99098133333333; N=35.972944; O=0.04490681078240458; P=0.04490681078240458; Q=331.5400618181244; R=2983.86055636312; S=0.01104158832504556; T=0.01104158832504556 },
    @{ A="M2"; B="Mif"; C="Cxcr4"; D="ECs"; E=3; F=1; G=24.21819066666667; H=72.654572; I=0.2153672905874018; J=0.2153672905874018; K=3; L=1; M=23.66111866666667; N=70.983356; O=0.08861204511346259; P=0.08861204511346259; Q=573.0294832559591; R=5157.265349303632; S=0.01908413606949506; T=0.01908413606949506 },
    @{ A="M2"; B="Mif"; C="Cxcr4"; D="FAPs"; E=3; F=1; G=24.21819066666667; H=72.654572; I=0.2153672905874018; J=0.2153672905874018; K=1; L=0.3333333333333333; M=0.1622346666666667; N=0.486704; O=0.0006075767508780888; P=0.0006075767508780888; Q=3.929030090076445; R=35.361270810688; S=0.0001308521586605108; T=0.0001308521586605108 },
    @{ A="M2"; B="Mif"; C="Cxcr4"; D="M2"; E=3; F=1; G=24.21819066666667; H=72.654572; I=0.2153672905874018; J=0.2153672905874018; K=3; L=1; M=231.2048796666667; N=693.614639; O=0.8658735673532548; P=0.8658735673532548; Q=5599.363858831057; R=50394.27472947951; S=0.1864808441921187; T=0.1864808441921187 },
    @{ A="M2"; B="Mif"; C="Cxcr4"; D="sCs"; E=3; F=1; G=24.21819066666667; H=72.654572; I=0.2153672905874018; J=0.2153672905874018; K=3; L=1; M=11.99098133333333; N=35.972944; O=0.04490681078240458; P=0.04490681078240458; Q=290.3998722111075; R=2613.598849899968; S=0.009671458167127597; T=0.009671458167127596 },
    @{ A="sCs"; B="Mif"; C="Cxcr4"; D="ECs"; E=3; F=1; G=36.47324633333334; H=109.419739; I=0.3243489305148018; J=0.3243489305148018; K=3; L=1; M=23.66111866666667; N=70.983356; O=0.08861204511346259; P=0.08861204511346259; Q=862.9978096515649; R=7766.980286864085; S=0.02874122206328096; T=0.02874122206328096 },
    @{ A="sCs"; B="Mif"; C="Cxcr4"; D="FAPs"; E=3; F=1; G=36.47324633333334; H=109.419739; I=0.3243489305148018; J=0.3243489305148018; K=1; L=0.3333333333333333; M=0.1622346666666667; N=0.486704; O=0.0006075767508780888; P=0.0006075767508780888; Q=5.917224961139556; R=53.255024650256; S=0.0001970668693529663; T=0.0001970668693529663 },
    @{ A="sCs"; B="Mif"; C="Cxcr4"; D="M2"; E=3; F=1; G=36.47324633333334; H=109.419739; I=0.3243489305148018; J=0.3243489305148018; K=3; L=1; M=231.2048796666667; N=693.614639; O=0.8658735673532548; P=0.8658735673532548; Q=8432.792529551025; R=75895.13276595922; S=0.2808451655320644; T=0.2808451655320644 },
    @{ A="sCs"; B="Mif"; C="Cxcr4"; D="sCs"; E=3; F=1; G=36.47324633333334; H=109.419739; I=0.3243489305148018; J=0.3243489305148018; K=3; L=1; M=11.99098133333333; N=35.972944; O=0.04490681078240458; P=0.04490681078240458; Q=437.3500159490684; R=3936.150143541616; S=0.0145654760501035; T=0.01456547605010349 }
)

$rowNum = 2
foreach ($rd in $rowsData) {
    $ws.Cells.Item($rowNum, 1).Value = $rd.A
    $ws.Cells.Item($rowNum, 2).Value = $rd.B
    $ws.Cells.Item($rowNum, 3).Value = $rd.C
    $ws.Cells.Item($rowNum, 4).Value = $rd.D
    $ws.Cells.Item($rowNum, 5).Value = $rd.E
    $ws.Cells.Item($rowNum, 6).Value = $rd.F
    $ws.Cells.Item($rowNum, 7).Value = $rd.G
    $ws.Cells.Item($rowNum, 8).Value = $rd.H
    $ws.Cells.Item($rowNum, 9).Value = $rd.I
    $ws.Cells.Item($rowNum, 10).Value = $rd.J
    $ws.Cells.Item($rowNum, 11).Value = $rd.K
    $ws.Cells.Item($rowNum, 12).Value = $rd.L
    $ws.Cells.Item($rowNum, 13).Value = $rd.M
    $ws.Cells.Item($rowNum, 14).Value = $rd.N
    $ws.Cells.Item($rowNum, 15).Value = $rd.O
    $ws.Cells.Item($rowNum, 16).Value = $rd.P
    $ws.Cells.Item($rowNum, 17).Value = $rd.Q
    $ws.Cells.Item($rowNum, 18).Value = $rd.R
    $ws.Cells.Item($rowNum, 19).Value = $rd.S
    $ws.Cells.Item($rowNum, 20).Value = $rd.T
    $rowNum = $rowNum + 1
}
